$wb = $excel.ActiveWorkbook

# --- Typography sheet: fill in the missing "Wildcard Characters" value for row 4 ---
$wsTypography = $wb.Worksheets.Item("Typography")
$wildcardChars = " !" + [char]8221 + [char]34 + "#*%&()'" + [char]36 + "+-@_, .:;?/~" + [char]177 + [char]215 + [char]247 + [char]8226 + [char]186 + [char]96 + [char]180 + "{}" + [char]169 + [char]163 + [char]8364 + "^" + [char]174 + [char]165 + "_=[]" + [char]161 + [char]162 + "|\" + [char]191 + "><"
$wsTypography.Range("G4").Value = $wildcardChars

# --- Translation sheet: add two new rows with translated UI text ---
$wsTranslation = $wb.Worksheets.Item("Translation")

$wsTranslation.Range("B75").Value = "SingleUseId81"
$wsTranslation.Range("C75").Value = "Typography_label"
$wsTranslation.Range("D75").Value = "Left"
$wsTranslation.Range("E75").Value = "LTR"
$wsTranslation.Range("F75").Value = "<value>"

$wsTranslation.Range("B76").Value = "SingleUseId83"
$wsTranslation.Range("C76").Value = "Typography_label"
$wsTranslation.Range("D76").Value = "Center"
$wsTranslation.Range("E76").Value = "LTR"
$wsTranslation.Range("F76").Value = "Ok"
